# Update distance/direction matrix values on the active sheet.
# Direction codes: 1 = Frente, 2 = Direita, 3 = Esquerda

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "AF6"  = 2
    "J12"  = 3
    "AT12" = 2
    "C13"  = 3
    "K13"  = 3
    "AS13" = 2
    "AN14" = 3
    "AO14" = 3
    "AQ14" = 3
    "AV14" = 3
    "I15"  = 3
    "M15"  = 3
    "Y15"  = 3
    "P35"  = 2
    "E36"  = 2
    "G40"  = 3
    "H40"  = 2
    "N40"  = 2
    "S40"  = 3
    "T40"  = 3
    "U40"  = 3
    "V40"  = 3
    "W40"  = 3
    "X40"  = 3
    "Z40"  = 2
    "AA40" = 2
    "AB40" = 2
    "AC40" = 2
    "AD40" = 2
    "AE40" = 2
    "AG40" = 2
    "AH40" = 3
    "AL40" = 2
    "AR40" = 2
    "AW40" = 3
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Preserve the sheet view's display options (gridlines stay visible).
$excel.ActiveWindow.DisplayGridlines = $true

# Update the view's scroll position (top-left visible cell) and the
# active cell selection to match the new location in the matrix.
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B14").Select()
